$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 10630
$ws.Range("I86").Value = 12365.1
$ws.Range("J86").Value = 7159.8
$ws.Range("K86").Value = 12365.1
$ws.Range("L86").Value = 7159.8
$ws.Range("M86").Value = -11242.1
$ws.Range("N86").Value = -9405.799999999999
$ws.Range("H89").Value = 10630
$ws.Range("I89").Value = 12365.1
$ws.Range("J89").Value = 7159.8
$ws.Range("K89").Value = 61825.5
$ws.Range("L89").Value = 35799
$ws.Range("M89").Value = -56209.5
$ws.Range("N89").Value = -47031
$ws.Range("H113").Value = 83335480
$ws.Range("I113").Value = 25002224
$ws.Range("J113").Value = 200002000
$ws.Range("K113").Value = 25002224
$ws.Range("L113").Value = 200002000
$ws.Range("M113").Value = -24998970
$ws.Range("N113").Value = -200008508
$ws.Range("H135").Value = 1247.4445
$ws.Range("I135").Value = 1291.4706
$ws.Range("J135").Value = 499
$ws.Range("K135").Value = 11623.2354
$ws.Range("L135").Value = 4491
$ws.Range("M135").Value = -9088.235400000001
$ws.Range("N135").Value = -9561
$ws.Range("H137").Value = 4257.074
$ws.Range("I137").Value = 1891
$ws.Range("J137").Value = 8279.4
$ws.Range("K137").Value = 5673
$ws.Range("L137").Value = 24838.2
$ws.Range("M137").Value = -3123
$ws.Range("N137").Value = -29938.2
$ws.Range("H138").Value = 2143.7358
$ws.Range("I138").Value = 1031.5883
$ws.Range("J138").Value = 2668.9167
$ws.Range("K138").Value = 3094.7649
$ws.Range("L138").Value = 8006.750100000001
$ws.Range("M138").Value = 2045.2351
$ws.Range("N138").Value = -18286.7501
$ws.Range("H141").Value = 3058.0908
$ws.Range("I141").Value = 3058.0908
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 9174.2724
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3994.2724

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1071.1818
$ws.Range("I4").Value = 920.55554
$ws.Range("J4").Value = 1749
$ws.Range("K4").Value = 920.55554
$ws.Range("L4").Value = 1749
$ws.Range("M4").Value = -804.55554
$ws.Range("N4").Value = -1981
$ws.Range("H5").Value = 363.75
$ws.Range("I5").Value = 257
$ws.Range("J5").Value = 598.6
$ws.Range("K5").Value = 257
$ws.Range("L5").Value = 598.6
$ws.Range("M5").Value = -145
$ws.Range("N5").Value = -822.6
$ws.Range("H32").Value = 15629677
$ws.Range("I32").Value = 17245100
$ws.Range("J32").Value = 13916
$ws.Range("K32").Value = 17245100
$ws.Range("L32").Value = 13916
$ws.Range("M32").Value = -17244813
$ws.Range("N32").Value = -14490
$ws.Range("H61").Value = 8950772
$ws.Range("I61").Value = 10207679
$ws.Range("J61").Value = 152422.42
$ws.Range("K61").Value = 10207679
$ws.Range("L61").Value = 152422.42
$ws.Range("M61").Value = -10207467
$ws.Range("N61").Value = -152846.42
$ws.Range("H108").Value = 92500
$ws.Range("I108").Value = 100000
$ws.Range("J108").Value = 70000
$ws.Range("K108").Value = 100000
$ws.Range("L108").Value = 70000
$ws.Range("M108").Value = -96160
$ws.Range("N108").Value = -77680
$ws.Range("H136").Value = 8950772
$ws.Range("I136").Value = 10207679
$ws.Range("J136").Value = 152422.42
$ws.Range("K136").Value = 30623037
$ws.Range("L136").Value = 457267.26
$ws.Range("M136").Value = -30620487
$ws.Range("N136").Value = -462367.26

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 363.75
$ws.Range("I4").Value = 257
$ws.Range("J4").Value = 598.6
$ws.Range("K4").Value = 257
$ws.Range("L4").Value = 598.6
$ws.Range("M4").Value = -142
$ws.Range("N4").Value = -828.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 33950
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 33950
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 33950
$ws.Range("N41").Value = -34806
$ws.Range("H50").Value = 67500
$ws.Range("I50").Value = 49000
$ws.Range("J50").Value = 73666.664
$ws.Range("K50").Value = 49000
$ws.Range("L50").Value = 73666.664
$ws.Range("M50").Value = -48375
$ws.Range("N50").Value = -74916.664
$ws.Range("H51").Value = 47499.832
$ws.Range("I51").Value = 24999.666
$ws.Range("J51").Value = 70000
$ws.Range("K51").Value = 24999.666
$ws.Range("L51").Value = 70000
$ws.Range("M51").Value = -24263.666
$ws.Range("N51").Value = -71472
$ws.Range("H58").Value = 1111.0555
$ws.Range("I58").Value = 1000.8
$ws.Range("J58").Value = 1248.875
$ws.Range("K58").Value = 1000.8
$ws.Range("L58").Value = 1248.875
$ws.Range("M58").Value = -797.8
$ws.Range("N58").Value = -1654.875
$ws.Range("H59").Value = 50000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 50000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 50000
$ws.Range("N59").Value = -52290
$ws.Range("H60").Value = 74000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 74000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 74000
$ws.Range("N60").Value = -75022
$ws.Range("H61").Value = 47499.832
$ws.Range("I61").Value = 24999.666
$ws.Range("J61").Value = 70000
$ws.Range("K61").Value = 24999.666
$ws.Range("L61").Value = 70000
$ws.Range("M61").Value = -24651.666
$ws.Range("N61").Value = -70696
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("H122").Value = 933.17645
$ws.Range("I122").Value = 1001.1818
$ws.Range("J122").Value = 808.5
$ws.Range("K122").Value = 3003.5454
$ws.Range("L122").Value = 2425.5
$ws.Range("M122").Value = -553.5454
$ws.Range("N122").Value = -7325.5
$ws.Range("H135").Value = 92142.71000000001
$ws.Range("I135").Value = 50000
$ws.Range("J135").Value = 99166.5
$ws.Range("K135").Value = 50000
$ws.Range("L135").Value = 99166.5
$ws.Range("M135").Value = -44930
$ws.Range("N135").Value = -109306.5
$ws.Range("H136").Value = 1111.0555
$ws.Range("I136").Value = 1000.8
$ws.Range("J136").Value = 1248.875
$ws.Range("K136").Value = 3002.4
$ws.Range("L136").Value = 3746.625
$ws.Range("M136").Value = -452.3999999999996
$ws.Range("N136").Value = -8846.625
$ws.Range("M60").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 936.25
$ws.Range("I109").Value = 803.6316
$ws.Range("J109").Value = 3456
$ws.Range("K109").Value = 2410.8948
$ws.Range("L109").Value = 10368
$ws.Range("M109").Value = -1370.8948
$ws.Range("N109").Value = -12448
$ws.Range("H131").Value = 12040.875
$ws.Range("I131").Value = 13732.25
$ws.Range("J131").Value = 10349.5
$ws.Range("K131").Value = 41196.75
$ws.Range("L131").Value = 31048.5
$ws.Range("M131").Value = -36156.75
$ws.Range("N131").Value = -41128.5
$ws.Range("H132").Value = 1365.909
$ws.Range("I132").Value = 914
$ws.Range("J132").Value = 1908.2
$ws.Range("K132").Value = 8226
$ws.Range("L132").Value = 17173.8
$ws.Range("M132").Value = -5696
$ws.Range("N132").Value = -22233.8
$ws.Range("H134").Value = 4607.273
$ws.Range("I134").Value = 3874.4119
$ws.Range("J134").Value = 7099
$ws.Range("K134").Value = 11623.2357
$ws.Range("L134").Value = 21297
$ws.Range("M134").Value = -6553.235700000001
$ws.Range("N134").Value = -31437
$ws.Range("H139").Value = 2207.6553
$ws.Range("I139").Value = 1327.75
$ws.Range("J139").Value = 2542.8572
$ws.Range("K139").Value = 3983.25
$ws.Range("L139").Value = 7628.571599999999
$ws.Range("M139").Value = 1156.75
$ws.Range("N139").Value = -17908.5716

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 21000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 21000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 21000
$ws.Range("N39").Value = -22064
$ws.Range("H92").Value = 25999.125
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 25999.125
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 25999.125
$ws.Range("N92").Value = -29743.125
$ws.Range("H109").Value = 47972.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 47972.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 47972.5
$ws.Range("N109").Value = -50052.5
$ws.Range("H132").Value = 142863710
$ws.Range("I132").Value = 142863710
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 428591130
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -428588600

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 49344.5
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 49344.5
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 49344.5
$ws.Range("N63").Value = -50842.5
$ws.Range("H66").Value = 49344.5
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 49344.5
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 148033.5
$ws.Range("N66").Value = -155521.5
$ws.Range("H132").Value = 345282.97
$ws.Range("I132").Value = 334958.16
$ws.Range("J132").Value = 500155
$ws.Range("K132").Value = 1004874.48
$ws.Range("L132").Value = 1500465
$ws.Range("M132").Value = -1002344.48
$ws.Range("N132").Value = -1505525
$ws.Range("H136").Value = 40032.195
$ws.Range("I136").Value = 4362.6875
$ws.Range("J136").Value = 78079.664
$ws.Range("K136").Value = 13088.0625
$ws.Range("L136").Value = 234238.992
$ws.Range("M136").Value = -10538.0625
$ws.Range("N136").Value = -239338.992
$ws.Range("H137").Value = 84000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 84000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 84000
$ws.Range("N137").Value = -94200

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 795
$ws.Range("I81").Value = 795
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1590
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -529
$ws.Range("H84").Value = 795
$ws.Range("I84").Value = 795
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 7950
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -2646
$ws.Range("H136").Value = 988.8
$ws.Range("I136").Value = 988.8
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2966.4
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -416.3999999999996
$ws.Range("H137").Value = 88000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 88000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 88000
$ws.Range("N137").Value = -98200
$ws.Range("N136").ClearContents()
